# feat(export): ajout d'une colonne pour indiquer les SMS
#
# Insert a new "Notifications SMS" column into the "Courriers" sheet
# (3rd worksheet), right before the existing "Courriers enregistrés"
# column, and make that sheet the active/selected tab.

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)

# Insert a new column G (shifting "Courriers enregistrés" and the rest
# one column to the right) and give it roughly the same width as the
# neighbouring column.
$ws3.Columns.Item(7).Insert()
$ws3.Columns.Item(7).ColumnWidth = $ws3.Columns.Item(6).ColumnWidth

# New header text for the inserted column.
$ws3.Cells.Item(2, 7).Value = "Notifications SMS"

# Make "Courriers" the active sheet/tab, with cell G3 selected.
$ws3.Activate()
$ws3.Range("G3").Select()
